# Update simulated transition-matrix probabilities after adding more
# simulated games / faster simulate-game logic (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1620689655172414
$ws.Range("C2").Value = 0.6241379310344828
$ws.Range("J2").Value = 0.01379310344827586
$ws.Range("P2").Value = 0.1241379310344828
$ws.Range("S2").Value = 0.07586206896551724
$ws.Range("B3").Value = 0.0053475935828877
$ws.Range("C3").Value = 0.03208556149732621
$ws.Range("J3").Value = 0.0267379679144385
$ws.Range("P3").Value = 0.7540106951871658
$ws.Range("S3").Value = 0.1818181818181818
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.6326530612244898
$ws.Range("S4").Value = 0.3469387755102041
$ws.Range("B6").Value = 0.05472636815920398
$ws.Range("D6").Value = 0.03482587064676617
$ws.Range("F6").Value = 0.09950248756218906
$ws.Range("J6").Value = 0.2338308457711443
$ws.Range("O6").Value = 0.03980099502487562
$ws.Range("Q6").Value = 0.1641791044776119
$ws.Range("R6").Value = 0.05970149253731343
$ws.Range("S6").Value = 0.3134328358208955
$ws.Range("B7").Value = 0.1325301204819277
$ws.Range("D7").Value = 0.02409638554216868
$ws.Range("E7").Value = 0.006024096385542169
$ws.Range("F7").Value = 0.03614457831325301
$ws.Range("J7").Value = 0.108433734939759
$ws.Range("O7").Value = 0.03012048192771084
$ws.Range("Q7").Value = 0.1506024096385542
$ws.Range("R7").Value = 0.09036144578313253
$ws.Range("S7").Value = 0.4216867469879518
$ws.Range("B8").Value = 0.1
$ws.Range("D8").Value = 0.00967741935483871
$ws.Range("F8").Value = 0.07419354838709677
$ws.Range("J8").Value = 0.1064516129032258
$ws.Range("O8").Value = 0.01290322580645161
$ws.Range("Q8").Value = 0.167741935483871
$ws.Range("R8").Value = 0.1129032258064516
$ws.Range("S8").Value = 0.4161290322580645
$ws.Range("B9").Value = 0.116504854368932
$ws.Range("D9").Value = 0.01456310679611651
$ws.Range("F9").Value = 0.05339805825242718
$ws.Range("J9").Value = 0.1359223300970874
$ws.Range("O9").Value = 0.01941747572815534
$ws.Range("Q9").Value = 0.1844660194174757
$ws.Range("R9").Value = 0.0970873786407767
$ws.Range("S9").Value = 0.3786407766990291
$ws.Range("B10").Value = 0.1356080489938758
$ws.Range("D10").Value = 0.02974628171478565
$ws.Range("E10").Value = 0.0008748906386701663
$ws.Range("F10").Value = 0.07436570428696412
$ws.Range("J10").Value = 0.1198600174978128
$ws.Range("O10").Value = 0.01487314085739283
$ws.Range("Q10").Value = 0.2178477690288714
$ws.Range("R10").Value = 0.07961504811898512
$ws.Range("S10").Value = 0.3272090988626422
$ws.Range("G11").Value = 0.1411290322580645
$ws.Range("J11").Value = 0.06451612903225806
$ws.Range("K11").Value = 0.1774193548387097
$ws.Range("L11").Value = 0.6048387096774194
$ws.Range("S11").Value = 0.01209677419354839
$ws.Range("G12").Value = 0.7133333333333334
$ws.Range("J12").Value = 0.2133333333333333
$ws.Range("K12").Value = 0.01333333333333333
$ws.Range("L12").Value = 0.01333333333333333
$ws.Range("S12").Value = 0.04666666666666667
$ws.Range("G13").Value = 0.7222222222222222
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.02777777777777778
$ws.Range("F15").Value = 0.03381642512077294
$ws.Range("H15").Value = 0.1207729468599034
$ws.Range("I15").Value = 0.07246376811594203
$ws.Range("J15").Value = 0.357487922705314
$ws.Range("K15").Value = 0.05797101449275362
$ws.Range("M15").Value = 0.02415458937198068
$ws.Range("O15").Value = 0.07246376811594203
$ws.Range("S15").Value = 0.2608695652173913
$ws.Range("F16").Value = 0.015
$ws.Range("H16").Value = 0.105
$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.405
$ws.Range("K16").Value = 0.12
$ws.Range("M16").Value = 0.02
$ws.Range("O16").Value = 0.06
$ws.Range("S16").Value = 0.15
$ws.Range("F17").Value = 0.007594936708860759
$ws.Range("H17").Value = 0.1493670886075949
$ws.Range("I17").Value = 0.1139240506329114
$ws.Range("J17").Value = 0.4683544303797468
$ws.Range("K17").Value = 0.08860759493670886
$ws.Range("M17").Value = 0.01265822784810127
$ws.Range("O17").Value = 0.06582278481012659
$ws.Range("S17").Value = 0.09367088607594937
$ws.Range("F18").Value = 0.01744186046511628
$ws.Range("H18").Value = 0.1569767441860465
$ws.Range("I18").Value = 0.1220930232558139
$ws.Range("J18").Value = 0.4069767441860465
$ws.Range("K18").Value = 0.09302325581395349
$ws.Range("M18").Value = 0.01744186046511628
$ws.Range("N18").Value = 0.005813953488372093
$ws.Range("O18").Value = 0.06395348837209303
$ws.Range("S18").Value = 0.1162790697674419
$ws.Range("F19").Value = 0.01792452830188679
$ws.Range("H19").Value = 0.1726415094339623
$ws.Range("I19").Value = 0.09150943396226414
$ws.Range("J19").Value = 0.389622641509434
$ws.Range("K19").Value = 0.1075471698113208
$ws.Range("M19").Value = 0.0160377358490566
$ws.Range("N19").Value = 0.0009433962264150943
$ws.Range("O19").Value = 0.07641509433962264
$ws.Range("S19").Value = 0.1273584905660377
